$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (after the weekly re-sort / update).
# Columns: D=Fecha, K=Variedad, L=Calidad, M=Volumen, N=Precio minimo,
#          O=Precio maximo, P=Precio promedio ponderado,
#          Q=Unidad de comercializacion, R=Origen, S=Precio $/Kg, T=Kg/unidad

$rows = @(
    @{ Row=2; D=44187; K="Dina";          L="Primera"; M=55; N=15000; O=16000; P=15455; Q="$/caja 15 kilos granel"; R="Región de O'Higgins"; S=1030; T=15 },
    @{ Row=3; D=44174; K="Castle Brite";  L="Primera"; M=75; N=9000;  O=10000; P=9467;  Q="$/caja 10 kilos";        R="Región de O'Higgins"; S=947;  T=10 },
    @{ Row=4; D=44176; K="Castle Brite";  L="Primera"; M=50; N=17000; O=18000; P=17400; Q="$/caja 18 kilos";        R="Región de O'Higgins"; S=967;  T=18 },
    @{ Row=5; D=44165; K="Castle Brite";  L="Segunda"; M=60; N=16000; O=17000; P=16500; Q="$/caja 15 kilos granel"; R="Provincia de Limarí";  S=1100; T=15 },
    @{ Row=6; D=44168; K="Castle Brite";  L="Primera"; M=30; N=16000; O=17000; P=16500; Q="$/caja 16 kilos granel"; R="Región de Coquimbo";   S=1031; T=16 },
    @{ Row=7; D=44181; K="Modesto";       L="Primera"; M=50; N=20000; O=21000; P=20500; Q="$/caja 18 kilos";        R="Región de Coquimbo";   S=1139; T=18 },
    @{ Row=8; D=44189; K="Dina";          L="Primera"; M=80; N=16000; O=17000; P=16562; Q="$/caja 18 kilos";        R="Región de O'Higgins"; S=920;  T=18 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($row, 11).Value = $r.K   # K - Variedad
    $ws.Cells.Item($row, 12).Value = $r.L   # L - Calidad
    $ws.Cells.Item($row, 13).Value = $r.M   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $r.N   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $r.O   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $r.P   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $r.R   # R - Origen
    $ws.Cells.Item($row, 19).Value = $r.S   # S - Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $r.T   # T - Kg / unidad
}
